$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

# New data for rows 2-11, columns B:G (ME, MAE, MSE, RMSE, SE, N)
# This reflects the ifoCAST full series evaluation: the error-statistics
# rows shift up by one period and a new final row (previously quarter Q9)
# is populated with fresh figures.
$data = @(
    @(-0.1106822175073481,  1.485348127483025,  8.723650662140232,  2.953582682462137,  3.020964860020825,  22),
    @(-0.5713011636365682,  1.053179700887051,  4.026630981043165,  2.006646700603563,  1.971105846864699,  21),
    @(-0.166036682263726,   0.6285186396184349, 0.8305896743800926, 0.911366926314584,  0.9193942482462926, 20),
    @(-0.1280209040135819,  0.7182792096296192, 0.930922228823257,  0.9648431109891685, 0.9825173072328391, 19),
    @(-0.1613809308236712,  0.6969901999678998, 0.7440183132057626, 0.8625649617308615, 0.8718992504770149, 18),
    @(-0.06310423029847312, 0.6271081964451903, 0.5990212621699853, 0.7739646388369338, 0.7951283312657454, 17),
    @(0.003107267299192557, 0.5393758537464697, 0.4187293164217235, 0.6470929735530463, 0.6683070442536863, 16),
    @(0.01079506215789655,  0.4821603736089785, 0.3882225284189086, 0.6230750584150425, 0.6448471538599968, 15),
    @(0.03588466164625848,  0.5517454314036715, 0.4354142101678629, 0.6598592351159926, 0.6837549615721854, 14),
    @(0.009683504591838199, 0.5061424767813365, 0.3656888178714484, 0.6047220997048548, 0.6293340143181445, 13)
)

$row = 2
foreach ($values in $data) {
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
    $row++
}
